# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" counts (column F) and, for two rows whose
# tickets have sold out, switches the "lowest price" column (G) from a
# numeric price to the text "不可售" (not available for sale).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 210
$ws.Cells.Item(3, 7).Value = "不可售"
$ws.Cells.Item(4, 6).Value = 95
$ws.Cells.Item(5, 6).Value = 1699
$ws.Cells.Item(6, 6).Value = 3283
$ws.Cells.Item(7, 6).Value = 901
$ws.Cells.Item(8, 6).Value = 2102
$ws.Cells.Item(9, 6).Value = 2013
$ws.Cells.Item(10, 6).Value = 1045
$ws.Cells.Item(11, 6).Value = 361
$ws.Cells.Item(13, 6).Value = 1628
$ws.Cells.Item(18, 6).Value = 104
$ws.Cells.Item(19, 6).Value = 1472
$ws.Cells.Item(20, 6).Value = 551
$ws.Cells.Item(22, 6).Value = 338
$ws.Cells.Item(23, 6).Value = 11851
$ws.Cells.Item(24, 6).Value = 11870
$ws.Cells.Item(25, 6).Value = 870
$ws.Cells.Item(26, 6).Value = 669
$ws.Cells.Item(27, 6).Value = 4
$ws.Cells.Item(28, 6).Value = 1864
$ws.Cells.Item(29, 6).Value = 165
$ws.Cells.Item(30, 6).Value = 475

# --- Sheet 2: 演出 (Performances) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 7

# --- Sheet 3: 本地生活 (Local Life) - no changes ---

# --- Sheet 4: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 210
$ws.Cells.Item(4, 7).Value = "不可售"
$ws.Cells.Item(6, 6).Value = 95
$ws.Cells.Item(7, 6).Value = 1699
$ws.Cells.Item(8, 6).Value = 3283
$ws.Cells.Item(9, 6).Value = 901
$ws.Cells.Item(10, 6).Value = 2102
$ws.Cells.Item(11, 6).Value = 2013
$ws.Cells.Item(12, 6).Value = 1045
$ws.Cells.Item(13, 6).Value = 361
$ws.Cells.Item(15, 6).Value = 1628
$ws.Cells.Item(19, 6).Value = 7
$ws.Cells.Item(22, 6).Value = 104
$ws.Cells.Item(23, 6).Value = 1472
$ws.Cells.Item(24, 6).Value = 551
$ws.Cells.Item(26, 6).Value = 338
$ws.Cells.Item(27, 6).Value = 11851
$ws.Cells.Item(28, 6).Value = 11870
$ws.Cells.Item(29, 6).Value = 870
$ws.Cells.Item(30, 6).Value = 669
$ws.Cells.Item(31, 6).Value = 4
$ws.Cells.Item(32, 6).Value = 1864
$ws.Cells.Item(35, 6).Value = 165
$ws.Cells.Item(36, 6).Value = 475
